$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Date (I) and Hours (K) columns to Text format so that values
# like "1/12/2021" or "12" are stored as literal strings, matching
# the original t="str" cell type, instead of being auto-converted
# to Excel dates/numbers.
$ws.Range("I2:I71").NumberFormat = "@"
$ws.Range("K2:K71").NumberFormat = "@"

$ws.Range("H2").Value = 1726.027397260274
$ws.Range("I2").Value = "28/11/2021"
$ws.Range("J2").Value = "Sunday"
$ws.Range("K2").Value = "12"

$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "28/11/2021"
$ws.Range("J3").Value = "Sunday"
$ws.Range("K3").Value = "0"

$ws.Range("H4").Value = 0
$ws.Range("I4").Value = "28/11/2021"
$ws.Range("J4").Value = "Sunday"
$ws.Range("K4").Value = "0"

$ws.Range("H5").Value = 37.442922374429216
$ws.Range("I5").Value = "28/11/2021"
$ws.Range("J5").Value = "Sunday"
$ws.Range("K5").Value = "2"

$ws.Range("H6").Value = 15.410958904109588
$ws.Range("I6").Value = "28/11/2021"
$ws.Range("J6").Value = "Sunday"
$ws.Range("K6").Value = "1"

$ws.Range("H7").Value = 42.12328767123288
$ws.Range("I7").Value = "28/11/2021"
$ws.Range("J7").Value = "Sunday"
$ws.Range("K7").Value = "3"

$ws.Range("H8").Value = 37.32876712328767
$ws.Range("I8").Value = "28/11/2021"
$ws.Range("J8").Value = "Sunday"
$ws.Range("K8").Value = "3"

$ws.Range("H9").Value = 6.84931506849315
$ws.Range("I9").Value = "28/11/2021"
$ws.Range("J9").Value = "Sunday"
$ws.Range("K9").Value = "1"

$ws.Range("H10").Value = 5.707762557077626
$ws.Range("I10").Value = "28/11/2021"
$ws.Range("J10").Value = "Sunday"
$ws.Range("K10").Value = "1"

$ws.Range("H11").Value = 5.479452054794519
$ws.Range("I11").Value = "28/11/2021"
$ws.Range("J11").Value = "Sunday"
$ws.Range("K11").Value = "1"

$ws.Range("H12").Value = 1869.86301369863
$ws.Range("I12").Value = "29/11/2021"
$ws.Range("J12").Value = "Monday"
$ws.Range("K12").Value = "13"

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "29/11/2021"
$ws.Range("J13").Value = "Monday"
$ws.Range("K13").Value = "0"

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "29/11/2021"
$ws.Range("J14").Value = "Monday"
$ws.Range("K14").Value = "0"

$ws.Range("H15").Value = 56.16438356164383
$ws.Range("I15").Value = "29/11/2021"
$ws.Range("J15").Value = "Monday"
$ws.Range("K15").Value = "3"

$ws.Range("H16").Value = 15.410958904109588
$ws.Range("I16").Value = "29/11/2021"
$ws.Range("J16").Value = "Monday"
$ws.Range("K16").Value = "1"

$ws.Range("H17").Value = 56.16438356164384
$ws.Range("I17").Value = "29/11/2021"
$ws.Range("J17").Value = "Monday"
$ws.Range("K17").Value = "4"

$ws.Range("H18").Value = 49.77168949771689
$ws.Range("I18").Value = "29/11/2021"
$ws.Range("J18").Value = "Monday"
$ws.Range("K18").Value = "4"

$ws.Range("H19").Value = 6.84931506849315
$ws.Range("I19").Value = "29/11/2021"
$ws.Range("J19").Value = "Monday"
$ws.Range("K19").Value = "1"

$ws.Range("H20").Value = 5.707762557077626
$ws.Range("I20").Value = "29/11/2021"
$ws.Range("J20").Value = "Monday"
$ws.Range("K20").Value = "1"

$ws.Range("H21").Value = 5.479452054794519
$ws.Range("I21").Value = "29/11/2021"
$ws.Range("J21").Value = "Monday"
$ws.Range("K21").Value = "1"

$ws.Range("H22").Value = 1869.86301369863
$ws.Range("I22").Value = "1/12/2021"
$ws.Range("J22").Value = "Tuesday"
$ws.Range("K22").Value = "13"

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = "1/12/2021"
$ws.Range("J23").Value = "Tuesday"
$ws.Range("K23").Value = "0"

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = "1/12/2021"
$ws.Range("J24").Value = "Tuesday"
$ws.Range("K24").Value = "0"

$ws.Range("H25").Value = 18.721461187214608
$ws.Range("I25").Value = "1/12/2021"
$ws.Range("J25").Value = "Tuesday"
$ws.Range("K25").Value = "1"

$ws.Range("H26").Value = 15.410958904109588
$ws.Range("I26").Value = "1/12/2021"
$ws.Range("J26").Value = "Tuesday"
$ws.Range("K26").Value = "1"

$ws.Range("H27").Value = 28.08219178082192
$ws.Range("I27").Value = "1/12/2021"
$ws.Range("J27").Value = "Tuesday"
$ws.Range("K27").Value = "2"

$ws.Range("H28").Value = 24.885844748858446
$ws.Range("I28").Value = "1/12/2021"
$ws.Range("J28").Value = "Tuesday"
$ws.Range("K28").Value = "2"

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = "1/12/2021"
$ws.Range("J29").Value = "Tuesday"
$ws.Range("K29").Value = ""

$ws.Range("H30").Value = 5.707762557077626
$ws.Range("I30").Value = "1/12/2021"
$ws.Range("J30").Value = "Tuesday"
$ws.Range("K30").Value = "1"

$ws.Range("H31").Value = 5.479452054794519
$ws.Range("I31").Value = "1/12/2021"
$ws.Range("J31").Value = "Tuesday"
$ws.Range("K31").Value = "1"

$ws.Range("H32").Value = 1438.3561643835617
$ws.Range("I32").Value = "2/12/2021"
$ws.Range("J32").Value = "Wednesday"
$ws.Range("K32").Value = "10"

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = "2/12/2021"
$ws.Range("J33").Value = "Wednesday"
$ws.Range("K33").Value = "0"

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = "2/12/2021"
$ws.Range("J34").Value = "Wednesday"
$ws.Range("K34").Value = "0"

$ws.Range("H35").Value = 37.442922374429216
$ws.Range("I35").Value = "2/12/2021"
$ws.Range("J35").Value = "Wednesday"
$ws.Range("K35").Value = "2"

$ws.Range("H36").Value = 15.410958904109588
$ws.Range("I36").Value = "2/12/2021"
$ws.Range("J36").Value = "Wednesday"
$ws.Range("K36").Value = "1"

$ws.Range("H37").Value = 28.08219178082192
$ws.Range("I37").Value = "2/12/2021"
$ws.Range("J37").Value = "Wednesday"
$ws.Range("K37").Value = "2"

$ws.Range("H38").Value = 24.885844748858446
$ws.Range("I38").Value = "2/12/2021"
$ws.Range("J38").Value = "Wednesday"
$ws.Range("K38").Value = "2"

$ws.Range("H39").Value = 13.6986301369863
$ws.Range("I39").Value = "2/12/2021"
$ws.Range("J39").Value = "Wednesday"
$ws.Range("K39").Value = "2"

$ws.Range("H40").Value = 5.707762557077626
$ws.Range("I40").Value = "2/12/2021"
$ws.Range("J40").Value = "Wednesday"
$ws.Range("K40").Value = "1"

$ws.Range("H41").Value = 5.479452054794519
$ws.Range("I41").Value = "2/12/2021"
$ws.Range("J41").Value = "Wednesday"
$ws.Range("K41").Value = "1"

$ws.Range("H42").Value = 1869.86301369863
$ws.Range("I42").Value = "3/12/2021"
$ws.Range("J42").Value = "Thursday"
$ws.Range("K42").Value = "13"

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = "3/12/2021"
$ws.Range("J43").Value = "Thursday"
$ws.Range("K43").Value = "0"

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = "3/12/2021"
$ws.Range("J44").Value = "Thursday"
$ws.Range("K44").Value = "0"

$ws.Range("H45").Value = 18.721461187214608
$ws.Range("I45").Value = "3/12/2021"
$ws.Range("J45").Value = "Thursday"
$ws.Range("K45").Value = "1"

$ws.Range("H46").Value = 15.410958904109588
$ws.Range("I46").Value = "3/12/2021"
$ws.Range("J46").Value = "Thursday"
$ws.Range("K46").Value = "1"

$ws.Range("H47").Value = 42.12328767123288
$ws.Range("I47").Value = "3/12/2021"
$ws.Range("J47").Value = "Thursday"
$ws.Range("K47").Value = "3"

$ws.Range("H48").Value = 37.32876712328767
$ws.Range("I48").Value = "3/12/2021"
$ws.Range("J48").Value = "Thursday"
$ws.Range("K48").Value = "3"

$ws.Range("H49").Value = 6.84931506849315
$ws.Range("I49").Value = "3/12/2021"
$ws.Range("J49").Value = "Thursday"
$ws.Range("K49").Value = "1"

$ws.Range("H50").Value = 5.707762557077626
$ws.Range("I50").Value = "3/12/2021"
$ws.Range("J50").Value = "Thursday"
$ws.Range("K50").Value = "1"

$ws.Range("H51").Value = 5.479452054794519
$ws.Range("I51").Value = "3/12/2021"
$ws.Range("J51").Value = "Thursday"
$ws.Range("K51").Value = "1"

$ws.Range("H52").Value = 1438.3561643835617
$ws.Range("I52").Value = "4/12/2021"
$ws.Range("J52").Value = "Friday"
$ws.Range("K52").Value = "10"

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = "4/12/2021"
$ws.Range("J53").Value = "Friday"
$ws.Range("K53").Value = "0"

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = "4/12/2021"
$ws.Range("J54").Value = "Friday"
$ws.Range("K54").Value = "0"

$ws.Range("H55").Value = 37.442922374429216
$ws.Range("I55").Value = "4/12/2021"
$ws.Range("J55").Value = "Friday"
$ws.Range("K55").Value = "2"

$ws.Range("H56").Value = 15.410958904109588
$ws.Range("I56").Value = "4/12/2021"
$ws.Range("J56").Value = "Friday"
$ws.Range("K56").Value = "1"

$ws.Range("H57").Value = 28.08219178082192
$ws.Range("I57").Value = "4/12/2021"
$ws.Range("J57").Value = "Friday"
$ws.Range("K57").Value = "2"

$ws.Range("H58").Value = 24.885844748858446
$ws.Range("I58").Value = "4/12/2021"
$ws.Range("J58").Value = "Friday"
$ws.Range("K58").Value = "2"

$ws.Range("H59").Value = 6.84931506849315
$ws.Range("I59").Value = "4/12/2021"
$ws.Range("J59").Value = "Friday"
$ws.Range("K59").Value = "1"

$ws.Range("H60").Value = 5.707762557077626
$ws.Range("I60").Value = "4/12/2021"
$ws.Range("J60").Value = "Friday"
$ws.Range("K60").Value = "1"

$ws.Range("H61").Value = 5.479452054794519
$ws.Range("I61").Value = "4/12/2021"
$ws.Range("J61").Value = "Friday"
$ws.Range("K61").Value = "1"

$ws.Range("H62").Value = 1869.86301369863
$ws.Range("I62").Value = "5/12/2021"
$ws.Range("J62").Value = "Saturday"
$ws.Range("K62").Value = "13"

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = "5/12/2021"
$ws.Range("J63").Value = "Saturday"
$ws.Range("K63").Value = "0"

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = "5/12/2021"
$ws.Range("J64").Value = "Saturday"
$ws.Range("K64").Value = "0"

$ws.Range("H65").Value = 18.721461187214608
$ws.Range("I65").Value = "5/12/2021"
$ws.Range("J65").Value = "Saturday"
$ws.Range("K65").Value = "1"

$ws.Range("H66").Value = 15.410958904109588
$ws.Range("I66").Value = "5/12/2021"
$ws.Range("J66").Value = "Saturday"
$ws.Range("K66").Value = "1"

$ws.Range("H67").Value = 56.16438356164384
$ws.Range("I67").Value = "5/12/2021"
$ws.Range("J67").Value = "Saturday"
$ws.Range("K67").Value = "4"

$ws.Range("H68").Value = 49.77168949771689
$ws.Range("I68").Value = "5/12/2021"
$ws.Range("J68").Value = "Saturday"
$ws.Range("K68").Value = "4"

$ws.Range("H69").Value = 6.84931506849315
$ws.Range("I69").Value = "5/12/2021"
$ws.Range("J69").Value = "Saturday"
$ws.Range("K69").Value = "1"

$ws.Range("H70").Value = 5.707762557077626
$ws.Range("I70").Value = "5/12/2021"
$ws.Range("J70").Value = "Saturday"
$ws.Range("K70").Value = "1"

$ws.Range("H71").Value = 5.479452054794519
$ws.Range("I71").Value = "5/12/2021"
$ws.Range("J71").Value = "Saturday"
$ws.Range("K71").Value = "1"
